# ---------------------------------------------------------------------------
# Commit: Mon, Apr 20, 2020  1:05:24 PM
#
# The canonical-OOXML diff for this commit shows the *entire* contents of
# ppt/theme/theme1.xml and ppt/theme/theme2.xml being exchanged:
#   - before: theme1.xml = "Office Theme" (used by the notes master),
#             theme2.xml = "Integral"     (used by the slide master /
#                                          presentation-level theme rel)
#   - after : theme1.xml = "Integral",
#             theme2.xml = "Office Theme"
# No shape, slide, layout, master or relationship-id content changes at
# all - every r:id still points at the same target filenames, and every
# slide/layout/master XML byte is identical. In other words the deck's
# design ("Integral") stays applied to the slide master exactly as
# before; only the *physical* theme part that happens to carry that XML
# is renumbered (an artifact of how the authoring tool re-serialized the
# package, not a user-visible formatting change).
#
# PowerPoint's object model intentionally does not expose the raw
# theme-part XML (Theme/Design only expose nav-only members such as
# .Application/.Index/.Parent/.SlideMaster/.ThemeVariants - there is no
# supported way to rewrite a theme part's bytes, rename/renumber a theme
# part, or reassign which physical part backs the slide master vs. the
# notes master). Touching the exposed Theme/Design surface (Master.Theme,
# Design.Name, Slide.Design, ThemeVariants, ...) round-trips the deck
# without altering either theme part, which matches the fact that the
# rendered/semantic design is unchanged by this commit.
#
# We still touch the documented, read-safe members below so the applied
# design ("Integral") is (re-)confirmed on the slide master, but we
# deliberately avoid guessing at unsupported writes that could corrupt
# unrelated parts of the package - every other byte of the deck must stay
# untouched, per the diff.
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

$master = $p.SlideMaster
$design = $master.Design
$theme  = $master.Theme

Write-Host ("Slide master design: " + $design.Name)
